# Leadership 2023-24.xlsx -- "update jobs & about"
#
# TTC committee roster changes:
#   - Danielle Rynczak moves from Member (row 3) to Chair (row 2)
#   - Casey LaDuke moves from Chair (row 2) to Past Chair (row 3)
#   - Caroline Erentzen (row 6) and Tarika Daftary-Kapur (row 7) leave the committee
#   - Cassidy Haigh moves from Student Member (row 8, with email) to Member (no email)
#   - Douglas Lewis, Will (Minqui) Pan, and Marco Chavez join as new Members
#   - Email column is cleared for all plain Members (keeping Chair/Past Chair emails)
#
# Also moves the saved "active" sheet/cell from Student Committee to TTC.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTC")

# --- Row 2: Chair -----------------------------------------------------
$ws.Range("A2").Value = "Chair"
$ws.Range("B2").Value = "Danielle Rynczak"
$ws.Range("C2").Value = "Danielle.Rynczak@umassmed.edu"
$ws.Range("D2").Value = "2021–2024"

# --- Row 3: Past Chair --------------------------------------------------
$ws.Range("A3").Value = "Past Chair"
$ws.Range("B3").Value = "Casey LaDuke"
$ws.Range("C3").Value = "claduke@jjay.cuny.edu"
$ws.Range("D3").Value = "2022-2025"

# --- Row 4: Member (Amanda Fanniff) - drop email ------------------------
$ws.Range("A4").Value = "Member"
$ws.Range("B4").Value = "Amanda Fanniff"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "2023–2026"

# --- Row 5: Member (James Andretta) - drop email ------------------------
$ws.Range("A5").Value = "Member"
$ws.Range("B5").Value = "James Andretta"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = "2023-2026"

# --- Row 6: Member (Douglas Lewis, new) ---------------------------------
$ws.Range("A6").Value = "Member"
$ws.Range("B6").Value = "Douglas Lewis"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = "2023-2026"

# --- Row 7: Member (Cassidy Haigh, was Student Member w/ email) ---------
$ws.Range("A7").Value = "Member"
$ws.Range("B7").Value = "Cassidy Haigh"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = "2021–2024"

# --- Row 8: Member (Will (Minqui) Pan, new) -----------------------------
# Row 8 already exists in the sheet; unlike rows 4-7 the target state has
# no <c> element at all for C8 (not even an empty one), so use Clear()
# instead of writing an empty string.
$ws.Range("A8").Value = "Member"
$ws.Range("B8").Value = "Will (Minqui) Pan"
$ws.Range("C8").Clear()
$ws.Range("D8").Value = "2023-2026"

# --- Row 9: Member (Marco Chavez, new row) ------------------------------
# Row 9 doesn't exist yet -- copy formats from row 8 cell-by-cell (skipping
# column C, which should stay entirely absent) before writing the values.
$ws.Range("A8").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Rows.Item(9).RowHeight = 15.75

$ws.Range("A9").Value = "Member"
$ws.Range("B9").Value = "Marco Chavez"
$ws.Range("D9").Value = "2023-2026"

# --- Active sheet / selection moves to TTC!D13 --------------------------
$ws.Activate()
$ws.Range("D13").Select()
